# Apply the "Added house elements for basic outline to calculator" edit.
#
# Summary of the change:
#  - On the "Pixel House Outline" sheet, rows 4-10 (column A) get real
#    section-description labels for the house outline elements (these had
#    either a placeholder "Sample" / "Enter your data into only the green
#    sections" text, or were blank).
#  - The pixel-type selector cell (A28) is changed from the "Brilliant Bulb"
#    item to the "Item #700 RGB Strip (30 LED/10 IC)" item, which cascades
#    through the dependent formulas automatically on recalculation.
#  - The active sheet/tab moves from "Intro" to "Pixel House Outline", with
#    the selection on that sheet moving to A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pixel House Outline")

# New house-outline section labels (rows 4-10, column A).
$ws.Range("A4").Value  = "Left Angle, Left Side"
$ws.Range("A5").Value  = "Left Angle, Right Side"
$ws.Range("A6").Value  = "Right Angle, Left Side"
$ws.Range("A7").Value  = "Right Angle, Right Side"
$ws.Range("A8").Value  = "Left Side of Garage"
$ws.Range("A9").Value  = "Top of Garage"
$ws.Range("A10").Value = "Right Side of Garage"

# Switch the selected pixel type used for the house outline.
$ws.Range("A28").Value = "Item #700 RGB Strip (30 LED/10 IC)"

# Make "Pixel House Outline" the active sheet/tab, with A11 selected.
$ws.Activate()
$ws.Range("A11").Select()
